$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 32 <- original row 33 (columns B:AD); column A (index) stays put
$ws.Range("B32").Value = 6859118
$ws.Range("C32").Value = "Denmark Division 3"
$ws.Range("D32").Value = 45156.58333333334
$ws.Range("E32").Value = "Vanlse"
$ws.Range("F32").Value = "Holbk"
$ws.Range("G32").Value = 1
$ws.Range("H32").Value = 1
$ws.Range("I32").Value = 1
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = "D"
$ws.Range("L32").Value = 1.4
$ws.Range("M32").Value = 4.5
$ws.Range("N32").Value = 6
$ws.Range("O32").Value = 1.45
$ws.Range("P32").Value = 4.5
$ws.Range("Q32").Value = 5.25
$ws.Range("R32").Value = -1.25
$ws.Range("S32").Value = 1.925
$ws.Range("T32").Value = 1.875
$ws.Range("U32").Value = 3.25
$ws.Range("V32").Value = 1.975
$ws.Range("W32").Value = 1.825
$ws.Range("X32").Value = -1
$ws.Range("Y32").Value = 3.5
$ws.Range("Z32").Value = -1
$ws.Range("AA32").Value = -1
$ws.Range("AB32").Value = 0.875
$ws.Range("AC32").Value = -1
$ws.Range("AD32").Value = 0.825

# Row 33 <- original row 32 (columns B:AD); column A (index) stays put
$ws.Range("B33").Value = 6858768
$ws.Range("C33").Value = "Denmark Division 3"
$ws.Range("D33").Value = 45156.58333333334
$ws.Range("E33").Value = "Holstebro"
$ws.Range("F33").Value = "Young Boys FD"
$ws.Range("G33").Value = 2
$ws.Range("H33").Value = 3
$ws.Range("I33").Value = 2
$ws.Range("J33").Value = 2
$ws.Range("K33").Value = "A"
$ws.Range("L33").Value = 4.333
$ws.Range("M33").Value = 3.8
$ws.Range("N33").Value = 1.65
$ws.Range("O33").Value = 4.5
$ws.Range("P33").Value = 4
$ws.Range("Q33").Value = 1.6
$ws.Range("R33").Value = 0.75
$ws.Range("S33").Value = 1.95
$ws.Range("T33").Value = 1.75
$ws.Range("U33").Value = 3
$ws.Range("V33").Value = 1.9
$ws.Range("W33").Value = 1.9
$ws.Range("X33").Value = -1
$ws.Range("Y33").Value = -1
$ws.Range("Z33").Value = 0.6000000000000001
$ws.Range("AA33").Value = -0.5
$ws.Range("AB33").Value = 0.375
$ws.Range("AC33").Value = 0.8999999999999999
$ws.Range("AD33").Value = -1

# Row 44 <- original row 45 (columns B:AD); column A (index) stays put
$ws.Range("B44").Value = 6859115
$ws.Range("C44").Value = "Denmark Division 3"
$ws.Range("D44").Value = 45171.33333333334
$ws.Range("E44").Value = "Holbk"
$ws.Range("F44").Value = "Young Boys FD"
$ws.Range("G44").Value = 1
$ws.Range("H44").Value = 1
$ws.Range("I44").Value = 1
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = "D"
$ws.Range("L44").Value = 5
$ws.Range("M44").Value = 4.333
$ws.Range("N44").Value = 1.5
$ws.Range("O44").Value = 4.2
$ws.Range("P44").Value = 3.8
$ws.Range("Q44").Value = 1.65
$ws.Range("R44").Value = 0.75
$ws.Range("S44").Value = 1.95
$ws.Range("T44").Value = 1.85
$ws.Range("U44").Value = 3
$ws.Range("V44").Value = 1.95
$ws.Range("W44").Value = 1.85
$ws.Range("X44").Value = -1
$ws.Range("Y44").Value = 2.8
$ws.Range("Z44").Value = -1
$ws.Range("AA44").Value = 0.95
$ws.Range("AB44").Value = -1
$ws.Range("AC44").Value = -1
$ws.Range("AD44").Value = 0.8500000000000001

# Row 45 <- original row 44 (columns B:AD); column A (index) stays put
$ws.Range("B45").Value = 6858777
$ws.Range("C45").Value = "Denmark Division 3"
$ws.Range("D45").Value = 45171.33333333334
$ws.Range("E45").Value = "Vanlse"
$ws.Range("F45").Value = "Ishoj"
$ws.Range("G45").Value = 2
$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 1
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = "H"
$ws.Range("L45").Value = 2.15
$ws.Range("M45").Value = 3.5
$ws.Range("N45").Value = 2.9
$ws.Range("O45").Value = 2.1
$ws.Range("P45").Value = 3.5
$ws.Range("Q45").Value = 2.9
$ws.Range("R45").Value = -0.25
$ws.Range("S45").Value = 1.9
$ws.Range("T45").Value = 1.9
$ws.Range("U45").Value = 3
$ws.Range("V45").Value = 1.825
$ws.Range("W45").Value = 1.975
$ws.Range("X45").Value = 1.1
$ws.Range("Y45").Value = -1
$ws.Range("Z45").Value = -1
$ws.Range("AA45").Value = 0.8999999999999999
$ws.Range("AB45").Value = -1
$ws.Range("AC45").Value = -1
$ws.Range("AD45").Value = 0.9750000000000001

# Row 62 <- original row 64 (columns B:AD); column A (index) stays put
$ws.Range("B62").Value = 6859109
$ws.Range("C62").Value = "Denmark Division 3"
$ws.Range("D62").Value = 45192.33333333334
$ws.Range("E62").Value = "Vejgaard B"
$ws.Range("F62").Value = "Avarta"
$ws.Range("G62").Value = 1
$ws.Range("H62").Value = 2
$ws.Range("I62").Value = 1
$ws.Range("J62").Value = 1
$ws.Range("K62").Value = "A"
$ws.Range("L62").Value = 1.75
$ws.Range("M62").Value = 3.6
$ws.Range("N62").Value = 4
$ws.Range("O62").Value = 1.727
$ws.Range("P62").Value = 3.6
$ws.Range("Q62").Value = 4
$ws.Range("R62").Value = -0.75
$ws.Range("S62").Value = 1.975
$ws.Range("T62").Value = 1.825
$ws.Range("U62").Value = 2.75
$ws.Range("V62").Value = 1.85
$ws.Range("W62").Value = 1.95
$ws.Range("X62").Value = -1
$ws.Range("Y62").Value = -1
$ws.Range("Z62").Value = 3
$ws.Range("AA62").Value = -1
$ws.Range("AB62").Value = 0.825
$ws.Range("AC62").Value = 0.425
$ws.Range("AD62").Value = -0.5

# Row 63 <- original row 62 (columns B:AD); column A (index) stays put
$ws.Range("B63").Value = 6858792
$ws.Range("C63").Value = "Denmark Division 3"
$ws.Range("D63").Value = 45192.33333333334
$ws.Range("E63").Value = "Ishoj"
$ws.Range("F63").Value = "Oure FA"
$ws.Range("G63").Value = 2
$ws.Range("H63").Value = 2
$ws.Range("I63").Value = 2
$ws.Range("J63").Value = 2
$ws.Range("K63").Value = "D"
$ws.Range("L63").Value = 1.45
$ws.Range("M63").Value = 4.5
$ws.Range("N63").Value = 5.5
$ws.Range("O63").Value = 1.45
$ws.Range("P63").Value = 4.5
$ws.Range("Q63").Value = 5.25
$ws.Range("R63").Value = -1.25
$ws.Range("S63").Value = 1.925
$ws.Range("T63").Value = 1.875
$ws.Range("U63").Value = 3.25
$ws.Range("V63").Value = 1.875
$ws.Range("W63").Value = 1.925
$ws.Range("X63").Value = -1
$ws.Range("Y63").Value = 3.5
$ws.Range("Z63").Value = -1
$ws.Range("AA63").Value = -1
$ws.Range("AB63").Value = 0.875
$ws.Range("AC63").Value = 0.875
$ws.Range("AD63").Value = -1

# Row 64 <- original row 63 (columns B:AD); column A (index) stays put
$ws.Range("B64").Value = 6858791
$ws.Range("C64").Value = "Denmark Division 3"
$ws.Range("D64").Value = 45192.33333333334
$ws.Range("E64").Value = "IF Lyseng"
$ws.Range("F64").Value = "Young Boys FD"
$ws.Range("G64").Value = 3
$ws.Range("H64").Value = 3
$ws.Range("I64").Value = 1
$ws.Range("J64").Value = 3
$ws.Range("K64").Value = "D"
$ws.Range("L64").Value = 2.8
$ws.Range("M64").Value = 3.6
$ws.Range("N64").Value = 2.15
$ws.Range("O64").Value = 2.9
$ws.Range("P64").Value = 3.4
$ws.Range("Q64").Value = 2.15
$ws.Range("R64").Value = 0.25
$ws.Range("S64").Value = 1.875
$ws.Range("T64").Value = 1.925
$ws.Range("U64").Value = 3
$ws.Range("V64").Value = 1.875
$ws.Range("W64").Value = 1.925
$ws.Range("X64").Value = -1
$ws.Range("Y64").Value = 2.4
$ws.Range("Z64").Value = -1
$ws.Range("AA64").Value = 0.4375
$ws.Range("AB64").Value = -0.5
$ws.Range("AC64").Value = 0.875
$ws.Range("AD64").Value = -1

# Row 190 <- original row 191 (columns B:AD); column A (index) stays put
$ws.Range("B190").Value = 8089057
$ws.Range("C190").Value = "Denmark Division 3"
$ws.Range("D190").Value = 45437.33333333334
$ws.Range("E190").Value = "Vanlse"
$ws.Range("F190").Value = "IF Lyseng"
$ws.Range("G190").Value = 3
$ws.Range("H190").Value = 2
$ws.Range("I190").Value = 1
$ws.Range("J190").Value = 1
$ws.Range("K190").Value = "H"
$ws.Range("L190").Value = 2.8
$ws.Range("M190").Value = 3.75
$ws.Range("N190").Value = 2.05
$ws.Range("O190").Value = 2.625
$ws.Range("P190").Value = 3.8
$ws.Range("Q190").Value = 2.2
$ws.Range("R190").Value = 0
$ws.Range("S190").Value = 2.075
$ws.Range("T190").Value = 1.725
$ws.Range("U190").Value = 3.25
$ws.Range("V190").Value = 1.85
$ws.Range("W190").Value = 1.95
$ws.Range("X190").Value = 1.625
$ws.Range("Y190").Value = -1
$ws.Range("Z190").Value = -1
$ws.Range("AA190").Value = 1.075
$ws.Range("AB190").Value = -1
$ws.Range("AC190").Value = 0.8500000000000001
$ws.Range("AD190").Value = -1

# Row 191 <- original row 190 (columns B:AD); column A (index) stays put
$ws.Range("B191").Value = 8089056
$ws.Range("C191").Value = "Denmark Division 3"
$ws.Range("D191").Value = 45437.33333333334
$ws.Range("E191").Value = "Vejgaard B"
$ws.Range("F191").Value = "Oure FA"
$ws.Range("G191").Value = 1
$ws.Range("H191").Value = 1
$ws.Range("I191").Value = 0
$ws.Range("J191").Value = 0
$ws.Range("K191").Value = "D"
$ws.Range("L191").Value = 1.8
$ws.Range("M191").Value = 3.75
$ws.Range("N191").Value = 3.5
$ws.Range("O191").Value = 1.6
$ws.Range("P191").Value = 4
$ws.Range("Q191").Value = 4.333
$ws.Range("R191").Value = -0.75
$ws.Range("S191").Value = 1.75
$ws.Range("T191").Value = 1.95
$ws.Range("U191").Value = 3.25
$ws.Range("V191").Value = 1.95
$ws.Range("W191").Value = 1.85
$ws.Range("X191").Value = -1
$ws.Range("Y191").Value = 3
$ws.Range("Z191").Value = -1
$ws.Range("AA191").Value = -1
$ws.Range("AB191").Value = 0.95
$ws.Range("AC191").Value = -1
$ws.Range("AD191").Value = 0.8500000000000001

# Row 196 <- original row 197 (columns B:AD); column A (index) stays put
$ws.Range("B196").Value = 8088916
$ws.Range("C196").Value = "Denmark Division 3"
$ws.Range("D196").Value = 45444.375
$ws.Range("E196").Value = "Frem"
$ws.Range("F196").Value = "VSK Arhus"
$ws.Range("G196").Value = 0
$ws.Range("H196").Value = 3
$ws.Range("K196").Value = "A"
$ws.Range("L196").Value = 1.75
$ws.Range("M196").Value = 3.6
$ws.Range("N196").Value = 4
$ws.Range("O196").Value = 1.75
$ws.Range("P196").Value = 3.7
$ws.Range("Q196").Value = 3.9
$ws.Range("R196").Value = -0.75
$ws.Range("S196").Value = 1.95
$ws.Range("T196").Value = 1.85
$ws.Range("U196").Value = 2.75
$ws.Range("V196").Value = 1.775
$ws.Range("W196").Value = 2.025
$ws.Range("X196").Value = -1
$ws.Range("Y196").Value = -1
$ws.Range("Z196").Value = 2.9
$ws.Range("AA196").Value = -1
$ws.Range("AB196").Value = 0.8500000000000001
$ws.Range("AC196").Value = 0.3875
$ws.Range("AD196").Value = -0.5

# Row 197 <- original row 196 (columns B:AD); column A (index) stays put
$ws.Range("B197").Value = 8089059
$ws.Range("C197").Value = "Denmark Division 3"
$ws.Range("D197").Value = 45444.375
$ws.Range("E197").Value = "Oure FA"
$ws.Range("F197").Value = "Vanlse"
$ws.Range("G197").Value = 1
$ws.Range("H197").Value = 1
$ws.Range("K197").Value = "D"
$ws.Range("L197").Value = 2.2
$ws.Range("M197").Value = 3.6
$ws.Range("N197").Value = 2.6
$ws.Range("O197").Value = 3.3
$ws.Range("P197").Value = 4.2
$ws.Range("Q197").Value = 1.75
$ws.Range("R197").Value = 0.75
$ws.Range("S197").Value = 1.825
$ws.Range("T197").Value = 1.975
$ws.Range("U197").Value = 3.25
$ws.Range("V197").Value = 1.95
$ws.Range("W197").Value = 1.85
$ws.Range("X197").Value = -1
$ws.Range("Y197").Value = 3.2
$ws.Range("Z197").Value = -1
$ws.Range("AA197").Value = 0.825
$ws.Range("AB197").Value = -1
$ws.Range("AC197").Value = -1
$ws.Range("AD197").Value = 0.8500000000000001
